$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 225, shifting the existing rows 225-347 down to 227-349.
$ws.Rows("225:226").Insert()

# Populate the first new row (225): Camote, "1a nueva(o)", origin Peru
$ws.Range("A225").Value2 = 8
$ws.Range("B225").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C225").Value2 = "Coquimbo"
$ws.Range("D225").Value2 = 44452
$ws.Range("E225").Value2 = 4
$ws.Range("F225").Value2 = 100112045
$ws.Range("G225").Value2 = "Zapallo"
$ws.Range("H225").Value2 = "Camote"
$ws.Range("I225").Value2 = "1a nueva(o)"
$ws.Range("J225").Value2 = 720
$ws.Range("K225").Value2 = 950
$ws.Range("L225").Value2 = 1000
$ws.Range("M225").Value2 = 975
$ws.Range("N225").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O225").Value2 = "Perú"
$ws.Range("P225").Value2 = 975
$ws.Range("Q225").Value2 = 1
$ws.Range("R225").Value2 = "Hortaliza"

# Populate the second new row (226): Camote, "2a nueva(o)", origin Peru
$ws.Range("A226").Value2 = 8
$ws.Range("B226").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C226").Value2 = "Coquimbo"
$ws.Range("D226").Value2 = 44452
$ws.Range("E226").Value2 = 4
$ws.Range("F226").Value2 = 100112045
$ws.Range("G226").Value2 = "Zapallo"
$ws.Range("H226").Value2 = "Camote"
$ws.Range("I226").Value2 = "2a nueva(o)"
$ws.Range("J226").Value2 = 520
$ws.Range("K226").Value2 = 850
$ws.Range("L226").Value2 = 900
$ws.Range("M226").Value2 = 875
$ws.Range("N226").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O226").Value2 = "Perú"
$ws.Range("P226").Value2 = 875
$ws.Range("Q226").Value2 = 1
$ws.Range("R226").Value2 = "Hortaliza"
